# Applies numeric corrections to the leve-profit calculation columns
# (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the
# scheduled market-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3567
$ws.Range("I43").Value = 6500
$ws.Range("J43").Value = 2100.5
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 2100.5
$ws.Range("M43").Value = -6431
$ws.Range("N43").Value = -2238.5
$ws.Range("H70").Value = 3023.5278
$ws.Range("I70").Value = 4835.294
$ws.Range("J70").Value = 1402.4736
$ws.Range("K70").Value = 14505.882
$ws.Range("L70").Value = 4207.4208
$ws.Range("M70").Value = -14235.882
$ws.Range("N70").Value = -4747.4208
$ws.Range("H73").Value = 3023.5278
$ws.Range("I73").Value = 4835.294
$ws.Range("J73").Value = 1402.4736
$ws.Range("K73").Value = 14505.882
$ws.Range("L73").Value = 4207.4208
$ws.Range("M73").Value = -13569.882
$ws.Range("N73").Value = -6079.4208
$ws.Range("H98").Value = 1114.6875
$ws.Range("I98").Value = 1089
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1089
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 409
$ws.Range("N98").Value = -4496
$ws.Range("H106").Value = 22225802
$ws.Range("I106").Value = 41667944
$ws.Range("J106").Value = 6214.2856
$ws.Range("K106").Value = 41667944
$ws.Range("L106").Value = 6214.2856
$ws.Range("M106").Value = -41667313
$ws.Range("N106").Value = -7476.2856
$ws.Range("H122").Value = 1114.6875
$ws.Range("I122").Value = 1089
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3267
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -817
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1392901.9
$ws.Range("I32").Value = 1606833.9
$ws.Range("J32").Value = 2343.625
$ws.Range("K32").Value = 1606833.9
$ws.Range("L32").Value = 2343.625
$ws.Range("M32").Value = -1606546.9
$ws.Range("N32").Value = -2917.625
$ws.Range("H45").Value = 1376.2858
$ws.Range("I45").Value = 1169.8182
$ws.Range("J45").Value = 2133.3333
$ws.Range("K45").Value = 1169.8182
$ws.Range("L45").Value = 2133.3333
$ws.Range("M45").Value = -792.8181999999999
$ws.Range("N45").Value = -2887.3333
$ws.Range("H74").Value = 184746.39
$ws.Range("I74").Value = 228041.14
$ws.Range("J74").Value = 65685.81
$ws.Range("K74").Value = 228041.14
$ws.Range("L74").Value = 65685.81
$ws.Range("M74").Value = -227167.14
$ws.Range("N74").Value = -67433.81
$ws.Range("H77").Value = 184746.39
$ws.Range("I77").Value = 228041.14
$ws.Range("J77").Value = 65685.81
$ws.Range("K77").Value = 1140205.7
$ws.Range("L77").Value = 328429.05
$ws.Range("M77").Value = -1135837.7
$ws.Range("N77").Value = -337165.05
$ws.Range("H110").Value = 1197.5714
$ws.Range("I110").Value = 1090.8788
$ws.Range("J110").Value = 1588.7778
$ws.Range("K110").Value = 1090.8788
$ws.Range("L110").Value = 1588.7778
$ws.Range("M110").Value = 954.1212
$ws.Range("N110").Value = -5678.7778
$ws.Range("H140").Value = 29500
$ws.Range("J140").Value = 29500
$ws.Range("L140").Value = 29500
$ws.Range("N140").Value = -39860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 24000
$ws.Range("J132").Value = 24000
$ws.Range("L132").Value = 24000
$ws.Range("N132").Value = -34120
$ws.Range("H140").Value = 59500
$ws.Range("J140").Value = 59500
$ws.Range("L140").Value = 59500
$ws.Range("N140").Value = -69860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3274.1282
$ws.Range("I31").Value = 2311.48
$ws.Range("K31").Value = 2311.48
$ws.Range("M31").Value = -2016.48
$ws.Range("H34").Value = 3274.1282
$ws.Range("I34").Value = 2311.48
$ws.Range("K34").Value = 2311.48
$ws.Range("M34").Value = -2109.48
$ws.Range("H99").Value = 61081.883
$ws.Range("I99").Value = 113222.664
$ws.Range("J99").Value = 2423.5
$ws.Range("K99").Value = 113222.664
$ws.Range("L99").Value = 2423.5
$ws.Range("M99").Value = -111724.664
$ws.Range("N99").Value = -5419.5
$ws.Range("H105").Value = 883.2308
$ws.Range("I105").Value = 852.8095
$ws.Range("J105").Value = 1011
$ws.Range("K105").Value = 852.8095
$ws.Range("L105").Value = 1011
$ws.Range("M105").Value = 894.1905
$ws.Range("N105").Value = -4505
$ws.Range("H126").Value = 61081.883
$ws.Range("I126").Value = 113222.664
$ws.Range("J126").Value = 2423.5
$ws.Range("K126").Value = 339667.992
$ws.Range("L126").Value = 7270.5
$ws.Range("M126").Value = -337197.992
$ws.Range("N126").Value = -12210.5
$ws.Range("H134").Value = 1518.7931
$ws.Range("I134").Value = 985.05554
$ws.Range("J134").Value = 2392.182
$ws.Range("K134").Value = 2955.16662
$ws.Range("L134").Value = 7176.545999999999
$ws.Range("M134").Value = -420.16662
$ws.Range("N134").Value = -12246.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1135
$ws.Range("I3").Value = 1135
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3405
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3293
$ws.Range("N3").ClearContents()
$ws.Range("H131").Value = 1388.54
$ws.Range("I131").Value = 1357
$ws.Range("J131").Value = 1396.425
$ws.Range("K131").Value = 4071
$ws.Range("L131").Value = 4189.275
$ws.Range("M131").Value = 969
$ws.Range("N131").Value = -14269.275
$ws.Range("H133").Value = 3239.3572
$ws.Range("I133").Value = 1255
$ws.Range("J133").Value = 4727.625
$ws.Range("K133").Value = 3765
$ws.Range("L133").Value = 14182.875
$ws.Range("M133").Value = 1295
$ws.Range("N133").Value = -24302.875
$ws.Range("H136").Value = 2103.0625
$ws.Range("I136").Value = 1340.8182
$ws.Range("J136").Value = 3780
$ws.Range("K136").Value = 4022.4546
$ws.Range("L136").Value = 11340
$ws.Range("M136").Value = 1077.5454
$ws.Range("N136").Value = -21540
$ws.Range("H137").Value = 4985.7646
$ws.Range("I137").Value = 864.64
$ws.Range("J137").Value = 16433.334
$ws.Range("K137").Value = 2593.92
$ws.Range("L137").Value = 49300.00199999999
$ws.Range("M137").Value = 2506.08
$ws.Range("N137").Value = -59500.00199999999
$ws.Range("H139").Value = 1895.7273
$ws.Range("I139").Value = 892.1667
$ws.Range("J139").Value = 3100
$ws.Range("K139").Value = 2676.5001
$ws.Range("L139").Value = 9300
$ws.Range("M139").Value = 2463.4999
$ws.Range("N139").Value = -19580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4628.457
$ws.Range("I80").Value = 5869.5
$ws.Range("J80").Value = 2973.7334
$ws.Range("K80").Value = 5869.5
$ws.Range("L80").Value = 2973.7334
$ws.Range("M80").Value = -4871.5
$ws.Range("N80").Value = -4969.7334
$ws.Range("H83").Value = 4628.457
$ws.Range("I83").Value = 5869.5
$ws.Range("J83").Value = 2973.7334
$ws.Range("K83").Value = 29347.5
$ws.Range("L83").Value = 14868.667
$ws.Range("M83").Value = -24355.5
$ws.Range("N83").Value = -24852.667
$ws.Range("H132").Value = 4162.9487
$ws.Range("I132").Value = 4755.722
$ws.Range("K132").Value = 14267.166
$ws.Range("M132").Value = -11737.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 36500
$ws.Range("J36").Value = 36500
$ws.Range("L36").Value = 36500
$ws.Range("N36").Value = -37624
$ws.Range("H40").Value = 2547.4546
$ws.Range("I40").Value = 2580
$ws.Range("J40").Value = 2222
$ws.Range("K40").Value = 2580
$ws.Range("L40").Value = 2222
$ws.Range("M40").Value = -2444
$ws.Range("N40").Value = -2494
$ws.Range("H46").Value = 2775.5557
$ws.Range("I46").Value = 1666.6666
$ws.Range("J46").Value = 4993.3335
$ws.Range("K46").Value = 1666.6666
$ws.Range("L46").Value = 4993.3335
$ws.Range("M46").Value = -1478.6666
$ws.Range("N46").Value = -5369.3335
$ws.Range("H68").Value = 3285.4348
$ws.Range("I68").Value = 2820.2
$ws.Range("J68").Value = 3643.3076
$ws.Range("K68").Value = 2820.2
$ws.Range("L68").Value = 3643.3076
$ws.Range("M68").Value = -2071.2
$ws.Range("N68").Value = -5141.3076
$ws.Range("H70").Value = 39000
$ws.Range("J70").Value = 39000
$ws.Range("L70").Value = 39000
$ws.Range("N70").Value = -39540
$ws.Range("H71").Value = 3285.4348
$ws.Range("I71").Value = 2820.2
$ws.Range("J71").Value = 3643.3076
$ws.Range("K71").Value = 14101
$ws.Range("L71").Value = 18216.538
$ws.Range("M71").Value = -10357
$ws.Range("N71").Value = -25704.538
$ws.Range("H73").Value = 39000
$ws.Range("J73").Value = 39000
$ws.Range("L73").Value = 39000
$ws.Range("N73").Value = -40872
$ws.Range("H122").Value = 1833.3334
$ws.Range("I122").Value = 1833.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5500.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3050.0002
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 20714.5
$ws.Range("J135").Value = 20714.5
$ws.Range("L135").Value = 20714.5
$ws.Range("N135").Value = -30854.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51476.332
$ws.Range("J46").Value = 51476.332
$ws.Range("L46").Value = 51476.332
$ws.Range("N46").Value = -51938.332
$ws.Range("H113").Value = 195.5
$ws.Range("I113").Value = 195.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 586.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1583.5
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 51476.332
$ws.Range("J134").Value = 51476.332
$ws.Range("L134").Value = 154428.996
$ws.Range("N134").Value = -159498.996
